# Laborator 29.11.2022 - pagini cu doar HTML si CSS care simuleaza aspectul
# proiectului final de FE
#
# - Fixes the typo "Sarb Maria" -> "Sirb Maria" in column A.
# - Marks attendance (TRUE) for "sapt 9" (column J) for the students that
#   attended that week: Bordas Norbert, Farcas Vasile, Fat Roberto,
#   Gavrilut Adriana, Rat Adrian, Sarb/Sirb Maria.
# - Leaves the active selection on J29, matching the final state of the
#   workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the name typo
$ws.Range("A30").Value = "Sirb Maria"

# Tick the "sapt 9" (column J) attendance checkbox for the relevant rows
$ws.Range("J3").Value = $true
$ws.Range("J7").Value = $true
$ws.Range("J8").Value = $true
$ws.Range("J10").Value = $true
$ws.Range("J28").Value = $true
$ws.Range("J30").Value = $true

# Update the active selection to match the saved workbook state
$ws.Range("J29").Select() | Out-Null
